# Add a 5th keyword/slug row to the sheet:
#   A5 = "passive income ideas"                              (brand new shared string)
#   B5 = "passive.income.nadi.myfirstdrawermenuproject"      (reuses the string already in B3)
#
# B3/B4 use a "wrap text" cell format (style index 1 in styles.xml) while the
# sheet's implicit default format (style index 0) is used elsewhere (e.g. the
# <cols> defaults). To reproduce that exactly we copy the *formatting only*
# from B4 into B5 before writing B5's value, and we leave A5 on the sheet's
# default format (so it keeps style index 0, same as the diff's s="0").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New keyword cell (A5) - plain/default formatting.
$ws.Range("A5").Value = "passive income ideas"

# Copy B4's formatting (font + wrap text) onto B5 so it reuses the existing
# cell style instead of Excel minting a brand new one, then set its text.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = "passive.income.nadi.myfirstdrawermenuproject"

# Rows 3 & 4 are taller (24pt, to fit the wrapped text); match that on row 5.
$ws.Rows.Item(5).RowHeight = 24

# Move the active selection to the newly added cell, like the source file.
$ws.Range("B5").Select()
